# Update gh-pages to output generated at 456a3b4
# Increment the "sales/views" figures in column F by 1 for the rows that
# changed in the refreshed data pull, across the affected worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F23").Value = 141
$ws1.Range("F25").Value = 7585
$ws1.Range("F44").Value = 311
$ws1.Range("F45").Value = 213
$ws1.Range("F48").Value = 128

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 176
$ws3.Range("F3").Value = 2561

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 176
$ws4.Range("F20").Value = 141
$ws4.Range("F23").Value = 7585
$ws4.Range("F42").Value = 311
$ws4.Range("F43").Value = 213
$ws4.Range("F46").Value = 128
